# Updated cryptos list row values (Coin/Link/Price/Volume(1h)) to match the
# latest scrape. Price (column D) values are forced to Text format so that
# numeric-looking strings such as "305.91" or "2.525.83" are written back as
# literal text, matching the inline-string cells already used in the sheet
# (column D mixes thousands-separated prices like "42.913.07" with plain
# decimals, so it must never be auto-coerced into a Double).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row => @{ column letter = new value }
$updates = [ordered]@{
    2 = @{ "D"="42.913.07"; "E"="  -0.56%  " }
    3 = @{ "D"="2.301.68"; "E"="  -0.27%  " }
    4 = @{ "E"="  +0.00%  " }
    5 = @{ "D"="305.91"; "E"="  +1.64%  " }
    6 = @{ "D"="97.33"; "E"="  -0.87%  " }
    7 = @{ "E"="  -1.77%  " }
    8 = @{ "E"="  +0.00%  " }
    9 = @{ "E"="  -2.64%  " }
    10 = @{ "D"="35.72"; "E"="  -0.66%  " }
    11 = @{ "D"="0.0798"; "E"="  +0.78%  " }
    12 = @{ "D"="18.19"; "E"="  +0.85%  " }
    13 = @{ "E"="  +1.03%  " }
    14 = @{ "D"="6.78"; "E"="  -1.50%  " }
    15 = @{ "D"="2.658.23"; "E"="  -0.38%  " }
    16 = @{ "D"="2.297.46"; "E"="  -3.14%  " }
    17 = @{ "D"="0.784"; "E"="  -0.86%  " }
    18 = @{ "D"="42.843.37"; "E"="  -0.51%  " }
    19 = @{ "D"="12.88"; "E"="  -2.31%  " }
    20 = @{ "D"="0.0₃0906"; "E"="  -0.53%  " }
    21 = @{ "E"="  -1.45%  " }
    22 = @{ "D"="67.65"; "E"="  -1.07%  " }
    23 = @{ "D"="236.85"; "E"="  -0.71%  " }
    24 = @{ "E"="  -1.63%  " }
    25 = @{ "D"="2.48"; "E"="  +2.35%  " }
    26 = @{ "E"="  -0.04%  " }
    27 = @{ "D"="4.03"; "E"="  +0.08%  " }
    28 = @{ "D"="25.44"; "E"="  +0.83%  " }
    29 = @{ "D"="166.56"; "E"="  -0.31%  " }
    30 = @{ "E"="  +1.09%  " }
    31 = @{ "E"="  -1.17%  " }
    32 = @{ "D"="32.98"; "E"="  -0.01%  " }
    33 = @{ "E"="  +0.10%  " }
    34 = @{ "E"="  +2.11%  " }
    35 = @{ "D"="5.01"; "E"="  -2.63%  " }
    36 = @{ "D"="17.37"; "E"="  -4.84%  " }
    37 = @{ "D"="2.40"; "E"="  -0.49%  " }
    38 = @{ "D"="0.0694"; "E"="  +0.59%  " }
    39 = @{ "D"="0.102"; "E"="  -0.85%  " }
    40 = @{ "E"="  -1.96%  " }
    41 = @{ "E"="  -1.35%  " }
    42 = @{ "E"="  -0.83%  " }
    43 = @{ "D"="2.012.60"; "E"="  +0.03%  " }
    44 = @{ "D"="0.0282"; "E"="  -2.03%  " }
    45 = @{ "B"="FraxShare"; "C"="https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"; "D"="10.01"; "E"="  -1.88%  " }
    46 = @{ "B"="EnergySwap"; "C"="https://coinranking.com/coin/SbWqqTui-+energyswap-ens"; "D"="17.97"; "E"="  +3.09%  " }
    47 = @{ "B"="NEARProtocol"; "C"="https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"; "D"="2.79"; "E"="  -2.13%  " }
    48 = @{ "B"="HuobiToken"; "C"="https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"; "D"="2.92"; "E"="  +5.22%  " }
    49 = @{ "B"="MultiversX"; "C"="https://coinranking.com/coin/omwkOTglq+multiversx-egld"; "D"="54.03"; "E"="  -0.81%  " }
    50 = @{ "B"="RocketPoolETH"; "C"="https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"; "D"="2.525.83"; "E"="  -0.32%  " }
    51 = @{ "B"="BitcoinSV"; "C"="https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"; "D"="72.10"; "E"="  -1.02%  " }
}

$colIndex = @{ "B" = 2; "C" = 3; "D" = 4; "E" = 5 }

foreach ($row in $updates.Keys) {
    foreach ($col in $updates[$row].Keys) {
        $cell = $ws.Cells.Item($row, $colIndex[$col])
        if ($col -eq "D") {
            # Price column: keep it textual, not a coerced number
            $cell.NumberFormat = "@"
        }
        $cell.Value = $updates[$row][$col]
    }
}
